# Responded to more posts.
# Set the "Actual time length to complete" (column C) for the
# "DQ1 response 2" task (row 9) on the week2 sheet to 20 minutes
# (expressed as a fraction of a day, matching the workbook's [h]:mm format).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week2")

$ws.Range("C9").Value = 20 / 1440

$wb.Save()
